$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SQL")

# --- New SQL reference text blocks (added rows 44-49) ---
$sqlSelectTcN = " select * from  OLE.portal_user p  join OLE.portal_user_tin pt`n on p.PORTAL_USER_ID=pt.PORTAL_USER_ID`n where p.TC_ACCEPT_IND='N' and p.STS_CD='A'  and pt.PROV_TIN_NBR not in ('010553448')`n and p.SSO_ID not in ('userTest808')`n order by `n p.CREAT_DTTM asc`n fetch first 1 rows only with ur"
$sqlUpdateTcN = " update ole.PORTAL_USER p`n set p.TC_ACCEPT_IND='N' where p.SSO_ID= '{`$id}'"
$sqlHavingTwo = "select pt.PROV_TIN_NBR, pt.ACCESS_LVL,count(*)  as totalUsers  , pu.EMAIL_ADR_TXT ,pu.FST_NM,pu.LST_NM`nfrom ole.ENROLLED_PROVIDER ep join  ole.PORTAL_USER_TIN pt   on ep.PROV_TIN_NBR=pt.PROV_TIN_NBR`njoin ole.PORTAL_USER pu on pu.PORTAL_USER_ID=pt.PORTAL_USER_ID`ngroup by pt.PROV_TIN_NBR,pt.ACCESS_LVL,pu.EMAIL_ADR_TXT,pu.FST_NM,pu.LST_NM,pu.STS_CD`nhaving count(*) =2  and pt.ACCESS_LVL='A'   and pu.STS_CD='A'`nfetch first row only`nwith ur"
$sqlEnrolledProvider = " select * from OLE.ENROLLED_PROVIDER p where p.ENRL_STS_CD='A' and  p.PROV_TIN_NBR not in ('{`$id}')`n order by p.CREAT_DTTM desc FETCH FIRST 1 ROW ONLY  "
$sqlHavingOne = "select pt.PROV_TIN_NBR, pt.ACCESS_LVL,count(*)  as totalUsers  , pu.EMAIL_ADR_TXT ,pu.FST_NM,pu.LST_NM`nfrom ole.ENROLLED_PROVIDER ep join  ole.PORTAL_USER_TIN pt   on ep.PROV_TIN_NBR=pt.PROV_TIN_NBR`njoin ole.PORTAL_USER pu on pu.PORTAL_USER_ID=pt.PORTAL_USER_ID`ngroup by pt.PROV_TIN_NBR,pt.ACCESS_LVL,pu.EMAIL_ADR_TXT,pu.FST_NM,pu.LST_NM,pu.STS_CD`nhaving count(*) =1 and pt.ACCESS_LVL='A'   and pu.STS_CD='A'`nfetch first row only`nwith ur"
$sqlUpdateTcY = " update ole.PORTAL_USER p`n set p.TC_ACCEPT_IND='Y' where p.SSO_ID= '{`$id}'"

# Populate cell values in the same order the original author typed them in,
# so new shared-string entries land in the expected order.
$ws.Range("A44").Value = "'43"
$ws.Range("C44").Value = "Getting user with terms and conditions indicator as N"
$ws.Range("B44").Value = $sqlSelectTcN

$ws.Range("B45").Value = $sqlUpdateTcN
$ws.Range("A45").Value = "'44"

$ws.Range("A46").Value = "'45"
$ws.Range("A47").Value = "'46"
$ws.Range("B47").Value = $sqlHavingTwo

$ws.Range("A48").Value = "'47"
$ws.Range("B48").Value = $sqlEnrolledProvider

$ws.Range("B46").Value = $sqlHavingOne

$ws.Range("B49").Value = $sqlUpdateTcY
$ws.Range("A49").Value = "'48"

# --- Formatting to match the existing "Sno"/"Query"/"Comments" columns ---
$ws.Range("B44:B49").WrapText = $true

$ws.Rows.Item(44).RowHeight = 100.8
$ws.Rows.Item(45).RowHeight = 28.8
$ws.Rows.Item(46).RowHeight = 144
$ws.Rows.Item(47).RowHeight = 144
$ws.Rows.Item(48).RowHeight = 43.2
$ws.Rows.Item(49).RowHeight = 28.8

# --- View state: scrolled down with A49 as the active cell ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 45
$ws.Range("A49").Select()
